# "unify the conception of DataNode, DataTable, Entity."
#
# The sheet in this workbook (currently called "Property") is being
# folded into the unified DataNode/DataTable/Entity naming scheme, so it
# is renamed to "DataNode". The rest of the upstream diff is just resave
# noise from opening the file in a different Excel build/locale
# (fileVersion/build numbers, absPath, window geometry, xr/xr2/xr3 GUID
# stamps, default font substitution, etc.) - nothing else about the
# actual data changes, so only the rename (plus the cursor ending up on
# a different cell, which is what that resave's sheetView selection
# shows) is reproduced here.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet: "Property" -> "DataNode"
$ws.Name = "DataNode"

# The resave also left the selection on D39 instead of A9.
[void]$ws.Range("D39").Select()
